# Update Excel data - 2024-11-22 06:26:37
# Refresh live crypto market data across all three sheets: per-asset
# price/marketcap/volume/24h-change on "Top 50 Cryptocurrencies", the
# mirrored market-cap figures on "Top 5 by Market Cap", and the derived
# summary stats text on "Summary".

$wb = $excel.ActiveWorkbook

# --- Sheet 1: Top 50 Cryptocurrencies -------------------------------------
$ws1 = $wb.Worksheets.Item("Top 50 Cryptocurrencies")

$ws1.Range("C2").Value = 98890 ; $ws1.Range("D2").Value = 1959698004531 ; $ws1.Range("E2").Value = 111556511395 ; $ws1.Range("F2").Value = 1.79611
$ws1.Range("C3").Value = 3368.92 ; $ws1.Range("D3").Value = 406338286340 ; $ws1.Range("E3").Value = 57069836564 ; $ws1.Range("F3").Value = 8.067830000000001
$ws1.Range("C4").Value = 1.001 ; $ws1.Range("D4").Value = 130970143446 ; $ws1.Range("E4").Value = 182938028684 ; $ws1.Range("F4").Value = 0.09354999999999999
$ws1.Range("C5").Value = 259.58 ; $ws1.Range("D5").Value = 123494602662 ; $ws1.Range("E5").Value = 14771296331 ; $ws1.Range("F5").Value = 8.85698
$ws1.Range("C6").Value = 630.27 ; $ws1.Range("D6").Value = 92034027631 ; $ws1.Range("E6").Value = 2506552932 ; $ws1.Range("F6").Value = 3.43071
$ws1.Range("D7").Value = 80278940768 ; $ws1.Range("E7").Value = 18503538652 ; $ws1.Range("F7").Value = 26.92868
$ws1.Range("C8").Value = 0.393778 ; $ws1.Range("D8").Value = 58027240960 ; $ws1.Range("E8").Value = 9748555380 ; $ws1.Range("F8").Value = 2.59959
$ws1.Range("C9").Value = 0.999411 ; $ws1.Range("D9").Value = 38359750699 ; $ws1.Range("E9").Value = 10847807727 ; $ws1.Range("F9").Value = 0.02544
$ws1.Range("C10").Value = 3368.18 ; $ws1.Range("D10").Value = 33037771601 ; $ws1.Range("E10").Value = 140321815 ; $ws1.Range("F10").Value = 7.94378
$ws1.Range("C11").Value = 0.884283 ; $ws1.Range("D11").Value = 31681007165 ; $ws1.Range("E11").Value = 3764850697 ; $ws1.Range("F11").Value = 12.03466
$ws1.Range("C12").Value = 0.199391 ; $ws1.Range("D12").Value = 17228471586 ; $ws1.Range("E12").Value = 1069521924 ; $ws1.Range("F12").Value = 0.93377
$ws1.Range("C13").Value = 36.25 ; $ws1.Range("D13").Value = 14850836971 ; $ws1.Range("E13").Value = 1041350077 ; $ws1.Range("F13").Value = 6.74635
$ws1.Range("C14").Value = 0.00002487 ; $ws1.Range("D14").Value = 14667878935 ; $ws1.Range("E14").Value = 1592407723 ; $ws1.Range("F14").Value = 3.34849
$ws1.Range("C15").Value = 3994.45 ; $ws1.Range("D15").Value = 14430087078 ; $ws1.Range("E15").Value = 170220205 ; $ws1.Range("F15").Value = 7.91847
$ws1.Range("C16").Value = 98404 ; $ws1.Range("D16").Value = 14398635396 ; $ws1.Range("E16").Value = 820480113 ; $ws1.Range("F16").Value = 1.45316
$ws1.Range("C17").Value = 5.53 ; $ws1.Range("D17").Value = 14095403038 ; $ws1.Range("E17").Value = 624264587 ; $ws1.Range("F17").Value = 2.79731
$ws1.Range("D18").Value = 10245919377 ; $ws1.Range("E18").Value = 2093423841 ; $ws1.Range("F18").Value = 0.87947
$ws1.Range("C19").Value = 495.73 ; $ws1.Range("D19").Value = 9835172618 ; $ws1.Range("E19").Value = 1749875615 ; $ws1.Range("F19").Value = -6.09258
$ws1.Range("C20").Value = 3363.77 ; $ws1.Range("D20").Value = 9656540862 ; $ws1.Range("E20").Value = 2085909929 ; $ws1.Range("F20").Value = 7.97274
$ws1.Range("C21").Value = 15.26 ; $ws1.Range("D21").Value = 9606079703 ; $ws1.Range("E21").Value = 1257251807 ; $ws1.Range("F21").Value = 4.59659
$ws1.Range("C22").Value = 0.0000213 ; $ws1.Range("D22").Value = 8960436990 ; $ws1.Range("E22").Value = 6723740508 ; $ws1.Range("F22").Value = 9.34882
$ws1.Range("D23").Value = 8956646967 ; $ws1.Range("E23").Value = 844240722 ; $ws1.Range("F23").Value = 9.40147
$ws1.Range("C24").Value = 0.287774 ; $ws1.Range("D24").Value = 8635355878 ; $ws1.Range("E24").Value = 2343798557 ; $ws1.Range("F24").Value = 21.1674
$ws1.Range("C25").Value = 8.81 ; $ws1.Range("D25").Value = 8139643571 ; $ws1.Range("E25").Value = 3400762 ; $ws1.Range("F25").Value = 3.85022
$ws1.Range("C26").Value = 5.78 ; $ws1.Range("D26").Value = 7046717485 ; $ws1.Range("E26").Value = 1003167791 ; $ws1.Range("F26").Value = 4.90486
$ws1.Range("C27").Value = 90.59 ; $ws1.Range("D27").Value = 6833848492 ; $ws1.Range("E27").Value = 1334339155 ; $ws1.Range("F27").Value = 2.55445
$ws1.Range("C28").Value = 12.04 ; $ws1.Range("D28").Value = 6423652211 ; $ws1.Range("E28").Value = 846792578 ; $ws1.Range("F28").Value = 3.46879
$ws1.Range("C29").Value = 3548.45 ; $ws1.Range("D29").Value = 6174101568 ; $ws1.Range("E29").Value = 105211848 ; $ws1.Range("F29").Value = 8.140940000000001
$ws1.Range("C30").Value = 9.33 ; $ws1.Range("D30").Value = 5610657297 ; $ws1.Range("E30").Value = 871047022 ; $ws1.Range("F30").Value = 5.54208
$ws1.Range("A31").Value = "Hedera" ; $ws1.Range("B31").Value = "hbar" ; $ws1.Range("C31").Value = 0.138456 ; $ws1.Range("D31").Value = 5260297033 ; $ws1.Range("E31").Value = 953925782 ; $ws1.Range("F31").Value = 10.34519
$ws1.Range("A32").Value = "USDS" ; $ws1.Range("B32").Value = "usds" ; $ws1.Range("C32").Value = 0.998479 ; $ws1.Range("D32").Value = 5235359178 ; $ws1.Range("E32").Value = 15996881 ; $ws1.Range("F32").Value = -0.59761
$ws1.Range("A33").Value = "Cronos" ; $ws1.Range("B33").Value = "cro" ; $ws1.Range("C33").Value = 0.189654 ; $ws1.Range("D33").Value = 5165132662 ; $ws1.Range("E33").Value = 147234964 ; $ws1.Range("F33").Value = 8.17057
$ws1.Range("C34").Value = 9.630000000000001 ; $ws1.Range("D34").Value = 4579534087 ; $ws1.Range("E34").Value = 273635929 ; $ws1.Range("F34").Value = 6.76845
$ws1.Range("C35").Value = 27.99 ; $ws1.Range("D35").Value = 4205246650 ; $ws1.Range("E35").Value = 864501403 ; $ws1.Range("F35").Value = 4.88331
$ws1.Range("C36").Value = 0.00005235 ; $ws1.Range("D36").Value = 3929105870 ; $ws1.Range("E36").Value = 1612558748 ; $ws1.Range("F36").Value = -0.81293
$ws1.Range("C37").Value = 0.152247 ; $ws1.Range("D37").Value = 3838980805 ; $ws1.Range("E37").Value = 150112818 ; $ws1.Range("F37").Value = 0.94309
$ws1.Range("C38").Value = 7.35 ; $ws1.Range("D38").Value = 3813718303 ; $ws1.Range("E38").Value = 430257724 ; $ws1.Range("F38").Value = 0.43552
$ws1.Range("C39").Value = 0.465096 ; $ws1.Range("D39").Value = 3708833312 ; $ws1.Range("E39").Value = 499273661 ; $ws1.Range("F39").Value = 6.25999
$ws1.Range("C40").Value = 501.34 ; $ws1.Range("D40").Value = 3702918639 ; $ws1.Range("E40").Value = 279730680 ; $ws1.Range("F40").Value = 3.40664
$ws1.Range("D41").Value = 3685203575 ; $ws1.Range("E41").Value = 223710665 ; $ws1.Range("F41").Value = -0.39308
$ws1.Range("C42").Value = 24.79 ; $ws1.Range("D42").Value = 3574017253 ; $ws1.Range("E42").Value = 31550147 ; $ws1.Range("F42").Value = 2.82256
$ws1.Range("C43").Value = 0.9993109999999999 ; $ws1.Range("D43").Value = 3444947266 ; $ws1.Range("E43").Value = 179507355 ; $ws1.Range("F43").Value = -0.05053
$ws1.Range("D44").Value = 3440542892 ; $ws1.Range("E44").Value = 304974175 ; $ws1.Range("F44").Value = 5.08724
$ws1.Range("C45").Value = 3.36 ; $ws1.Range("D45").Value = 3353709873 ; $ws1.Range("E45").Value = 1288616085 ; $ws1.Range("F45").Value = 5.23119
$ws1.Range("D46").Value = 3341022642 ; $ws1.Range("E46").Value = 483200976 ; $ws1.Range("F46").Value = 2.98349
$ws1.Range("C47").Value = 0.786529 ; $ws1.Range("D47").Value = 3225954936 ; $ws1.Range("E47").Value = 1673163724 ; $ws1.Range("F47").Value = 13.77642
$ws1.Range("C48").Value = 160.51 ; $ws1.Range("D48").Value = 2968018563 ; $ws1.Range("E48").Value = 85853749 ; $ws1.Range("F48").Value = -0.55401
$ws1.Range("D49").Value = 2944746477 ; $ws1.Range("E49").Value = 350307514 ; $ws1.Range("F49").Value = 2.22298
$ws1.Range("A50").Value = "Filecoin" ; $ws1.Range("B50").Value = "fil" ; $ws1.Range("C50").Value = 4.68 ; $ws1.Range("D50").Value = 2816547654 ; $ws1.Range("E50").Value = 574134364 ; $ws1.Range("F50").Value = 6.59424
$ws1.Range("A51").Value = "Mantle" ; $ws1.Range("B51").Value = "mnt" ; $ws1.Range("C51").Value = 0.835849 ; $ws1.Range("D51").Value = 2813809441 ; $ws1.Range("E51").Value = 186127050 ; $ws1.Range("F51").Value = 14.36165

# --- Sheet 2: Top 5 by Market Cap ------------------------------------------
$ws2 = $wb.Worksheets.Item("Top 5 by Market Cap")
$ws2.Range("B2").Value = 1959698004531
$ws2.Range("B3").Value = 406338286340
$ws2.Range("B4").Value = 130970143446
$ws2.Range("B5").Value = 123494602662
$ws2.Range("B6").Value = 92034027631

# --- Sheet 3: Summary --------------------------------------------------
# These cells hold plain text ("$4345.48", "XRP (26.93%)", ...). Excel's
# COM value-coercion treats a leading "$" as a currency literal and would
# otherwise silently convert the cell to a formatted number, so the
# number format is forced to text ("@") for the write and then restored
# to the sheet's normal (General) style afterwards.
$ws3 = $wb.Worksheets.Item("Summary")

$ws3.Range("B2").NumberFormat = "@"
$ws3.Range("B2").Value = "`$4345.48"
$ws3.Range("B2").Style = "Normal"

$ws3.Range("B3").Value = "XRP (26.93%)"
$ws3.Range("B4").Value = "Bitcoin Cash (-6.09%)"
